$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Laatst bijgewerkt" timestamp (row 2, col A) ---
$ws.Range("A2").Value = "Laatst bijgewerkt: 2025-09-06 19:17:39"

# --- Remove existing hyperlinks so they can be rebuilt cleanly for the
#     shifted rows (a new match row is inserted above the old data) ---
$ws.Range("A1:N20").Hyperlinks.Delete()

# The "Beste inzet" (K) column holds plain text like "€4.88". Pre-format
# those cells as Text so Excel's smart-entry doesn't convert them into a
# currency number; the format gets reset back to Normal right after the
# values are written.
foreach ($r in 3..7) {
    $ws.Range("K$r").NumberFormat = "@"
}

# --- Row 3: new match "Duitsland vs Noord-Ierland" (inserted above the
#     previously existing rows, which all shift down by one) ---
$ws.Range("A3").Value = "Duitsland vs Noord-Ierland"
$ws.Range("B3").Value = "totaal aantal schoten op doel"
$ws.Range("C3").Value = "wedstrijd"
$ws.Range("D3").Value = "meer dan 10.5"
$ws.Range("E3").Value = "toto"
$ws.Range("F3").Value = 2.5
$ws.Range("G3").Value = "minder dan 10.5"
$ws.Range("H3").Value = "starcasino"
$ws.Range("I3").Value = 1.76
$ws.Range("J3").Value = "1=62, 2=88"
$ws.Range("K3").Value = "€4.88"
$ws.Range("L3").Value = 3.18
$ws.Range("M3").Value = "https://sport.toto.nl/wedden/wedstrijd/8778584"
$ws.Range("N3").Value = "https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394"

# --- Row 4: formerly row 3 (Litouwen vs Nederland / sem steijn) ---
$ws.Range("A4").Value = "Litouwen vs Nederland"
$ws.Range("B4").Value = "totaal aantal schoten"
$ws.Range("C4").Value = "sem steijn"
$ws.Range("D4").Value = "meer dan 3.5"
$ws.Range("E4").Value = "vbet"
$ws.Range("F4").Value = 1.91
$ws.Range("G4").Value = "minder dan 3.5"
$ws.Range("H4").Value = "jacks"
$ws.Range("I4").Value = 2.23
$ws.Range("J4").Value = "1=81, 2=69"
$ws.Range("K4").Value = "€3.87"
$ws.Range("L4").Value = 2.8
$ws.Range("M4").Value = "https://www.vbet.nl/nl/sports/pre-match/event-view/Soccer/World/18277589/world-cup-europe-qualification/27857408/litouwen-nederland"
$ws.Range("N4").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

# --- Row 5: formerly row 4 (Litouwen vs Nederland / nederland, 8.5) ---
$ws.Range("A5").Value = "Litouwen vs Nederland"
$ws.Range("B5").Value = "totaal aantal schoten op doel"
$ws.Range("C5").Value = "nederland"
$ws.Range("D5").Value = "meer dan 8.5"
$ws.Range("E5").Value = "toto"
$ws.Range("F5").Value = 2.45
$ws.Range("G5").Value = "minder dan 8.5"
$ws.Range("H5").Value = "jacks"
$ws.Range("I5").Value = 1.74
$ws.Range("J5").Value = "1=62, 2=88"
$ws.Range("K5").Value = "€1.9"
$ws.Range("L5").Value = 1.71
$ws.Range("M5").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N5").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

# --- Row 6: formerly row 5 (Litouwen vs Nederland / nederland, 7.5) ---
$ws.Range("A6").Value = "Litouwen vs Nederland"
$ws.Range("B6").Value = "totaal aantal schoten op doel"
$ws.Range("C6").Value = "nederland"
$ws.Range("D6").Value = "meer dan 7.5"
$ws.Range("E6").Value = "toto"
$ws.Range("F6").Value = 1.85
$ws.Range("G6").Value = "minder dan 7.5"
$ws.Range("H6").Value = "jacks"
$ws.Range("I6").Value = 2.25
$ws.Range("J6").Value = "1=82, 2=68"
$ws.Range("K6").Value = "€1.7"
$ws.Range("L6").Value = 1.5
$ws.Range("M6").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N6").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

# --- Row 7: formerly row 6 (Litouwen vs Nederland / nederland, 9.5) ---
$ws.Range("A7").Value = "Litouwen vs Nederland"
$ws.Range("B7").Value = "totaal aantal schoten op doel"
$ws.Range("C7").Value = "nederland"
$ws.Range("D7").Value = "meer dan 9.5"
$ws.Range("E7").Value = "toto"
$ws.Range("F7").Value = 3.5
$ws.Range("G7").Value = "minder dan 9.5"
$ws.Range("H7").Value = "jacks"
$ws.Range("I7").Value = 1.43
$ws.Range("J7").Value = "1=44, 2=106"
$ws.Range("K7").Value = "€1.58"
$ws.Range("L7").Value = 1.5
$ws.Range("M7").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N7").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

# Put the K column back to its normal (default) style now that the text
# values are safely stored as literal text.
foreach ($r in 3..7) {
    $ws.Range("K$r").Style = "Normal"
}

# --- Rebuild hyperlinks for M/N columns on rows 3-7 in order, so rIds
#     come out sequentially matching the canonical save order ---
$ws.Hyperlinks.Add($ws.Range("M3"), "https://sport.toto.nl/wedden/wedstrijd/8778584")
$ws.Hyperlinks.Add($ws.Range("N3"), "https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394")

$ws.Hyperlinks.Add($ws.Range("M4"), "https://www.vbet.nl/nl/sports/pre-match/event-view/Soccer/World/18277589/world-cup-europe-qualification/27857408/litouwen-nederland")
$ws.Hyperlinks.Add($ws.Range("N4"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")

$ws.Hyperlinks.Add($ws.Range("M5"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N5"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")

$ws.Hyperlinks.Add($ws.Range("M6"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N6"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")

$ws.Hyperlinks.Add($ws.Range("M7"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N7"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")

# Re-apply the built-in "Hyperlink" cell style to the link columns (style
# that was already used by the pre-existing link cells) now that the
# Hyperlinks.Add calls have run.
$ws.Range("M3:N7").Style = "Hyperlink"
